$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: Replace the full country-name lists in columns P and R (both the
# "Table 1 updated" and "Table 1 complete" sheets) with ISO3166-1 alpha-3
# country codes. Row 1 (header) is left untouched.
# ---------------------------------------------------------------------------
$sheetNames = @("Table 1 updated", "Table 1 complete")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range('P2').Value = 'AUS, CAN, CHL, CRI, DNK, FRA, NZL, NOR, KOR, ZAF, GBR, USA'; $ws.Range('R2').Value = 'AUS, CAN, CHL, CRI, DNK, FRA, NZL, NOR, KOR, ZAF, GBR, USA'
    $ws.Range('P3').Value = 'CUB, MEX, DOM, PRI, CHN, PER, IND, VEN, AUT, AUS, GRC, ESP, FIN, POL, GBR, RUS, LTU, CZE, USA, JPN, GHA, NLD, KOR, ZAF, DNK, HUN, SWE, CHE, IRL, EST, BEL, ISR, SVN, HRV, LUX, PRT, FRA, DEU, ITA'; $ws.Range('R3').Value = 'CUB, MEX, DOM, PRI, CHN, PER, IND, VEN, AUT, AUS, GRC, ESP, FIN, POL, GBR, RUS, LTU, CZE, USA, JPN, GHA, NLD, KOR, ZAF, DNK, HUN, SWE, CHE, IRL, EST, BEL, ISR, SVN, HRV, LUX, PRT, FRA, DEU, ITA'
    $ws.Range('P4').Value = 'NLD'; $ws.Range('R4').Value = 'NLD'
    $ws.Range('P5').Value = 'CAN, DEU, ITA, IRL, GBR, FIN, NLD, NOR'; $ws.Range('R5').Value = 'CAN, DEU, ITA, IRL, GBR, FIN, NLD, NOR'
    $ws.Range('P6').Value = 'CAN'; $ws.Range('R6').Value = 'CAN'
    $ws.Range('P7').Value = 'SWE'; $ws.Range('R7').Value = 'SWE'
    $ws.Range('P8').Value = 'ESP'; $ws.Range('R8').Value = 'ESP'
    $ws.Range('P9').Value = 'GBR, FIN'; $ws.Range('R9').Value = 'GBR, FIN'
    $ws.Range('P10').Value = 'NLD, ESP, GBR, DNK, SWE, BEL, NOR, FRA, FIN, GRC, DEU, CYP'; $ws.Range('R10').Value = 'NLD, ESP, GBR, DNK, SWE, BEL, NOR, FRA, FIN, GRC, DEU, CYP'
    $ws.Range('P11').Value = 'DEU, ITA, NLD, ESP, GBR, SWE'; $ws.Range('R11').Value = 'DEU, ITA, NLD, ESP, GBR, SWE'
    $ws.Range('P12').Value = 'CAN, NOR, FRA, DNK, EST, DEU, GRC, ITA, LTU, POL, RUS, IRL, PRT, BEL, ESP, SWE, NLD, GBR, FIN, AUS'; $ws.Range('R12').Value = 'CAN, NOR, FRA, DNK, EST, DEU, GRC, ITA, LTU, POL, RUS, IRL, PRT, BEL, ESP, SWE, NLD, GBR, FIN, AUS'
    $ws.Range('P13').Value = 'DEU, GBR, FIN, EST, DNK, ITA, ESP, LTU, AUS, SWE, RUS, POL, FRA, NOR, CAN'; $ws.Range('R13').Value = 'DEU, GBR, FIN, EST, DNK, ITA, ESP, LTU, AUS, SWE, RUS, POL, FRA, NOR, CAN'
    $ws.Range('P14').Value = 'USA, MEX, ISR, CRI, KOR, JPN, CHN, IND, GBR, AUT, BEL, HRV, CYP, CZE, DNK, EST, FIN, FRA, DEU, GRC, HUN, ITA, LVA, LTU, LUX, MLT, NLD, POL, PRT, ROU, SVK, SVN, ESP, SWE, CHE, IRL, IDN, GHA, RUS, ZAF, BRA, THA, MYS'; $ws.Range('R14').Value = 'USA, MEX, ISR, CRI, KOR, JPN, CHN, IND, GBR, AUT, BEL, HRV, CYP, CZE, DNK, EST, FIN, FRA, DEU, GRC, HUN, ITA, LVA, LTU, LUX, MLT, NLD, POL, PRT, ROU, SVK, SVN, ESP, SWE, CHE, IRL, IDN, GHA, RUS, ZAF, BRA, THA, MYS'
    $ws.Range('P15').Value = 'DNK, FIN, ITA, NLD, NOR, GBR, AUS, SWE'; $ws.Range('R15').Value = 'DNK, FIN, ITA, NLD, NOR, GBR, AUS, SWE'
    $ws.Range('P16').Value = 'GBR'; $ws.Range('R16').Value = 'GBR'
    $ws.Range('P17').Value = 'FRA, GRC, LTU, NOR, ESP, GBR'; $ws.Range('R17').Value = 'FRA, GRC, LTU, NOR, ESP, GBR'
    $ws.Range('P18').Value = 'CAN, CHN, IND, ZAF'; $ws.Range('R18').Value = 'CAN, CHN, IND, ZAF'
    $ws.Range('P19').Value = 'USA, FIN, AUS'; $ws.Range('R19').Value = 'USA, FIN, AUS'
    $ws.Range('P20').Value = 'AUS, DNK, BRA, CHN, FRA, ITA, ISR, JPN, NOR, GBR, USA'; $ws.Range('R20').Value = 'AUS, DNK, BRA, CHN, FRA, ITA, ISR, JPN, NOR, GBR, USA'
    $ws.Range('P22').Value = 'SWE, GBR, FIN, FRA, DNK'; $ws.Range('R22').Value = 'SWE, GBR, FIN, FRA, DNK'
    $ws.Range('P23').Value = 'SWE, DNK, FIN, USA, AUS'; $ws.Range('R23').Value = 'SWE, DNK, FIN, USA, AUS'
    $ws.Range('P26').Value = 'AUS, GBR, DNK, FRA, JPN, USA, LBN, ESP, SWE, MAR, CHN'; $ws.Range('R26').Value = 'AUS, GBR, DNK, FRA, JPN, USA, LBN, ESP, SWE, MAR, CHN'
    $ws.Range('P27').Value = 'GBR, IRL, PRT, FRA, CHE, ITA, FIN, USA, AUS'; $ws.Range('R27').Value = 'GBR, IRL, PRT, FRA, CHE, ITA, FIN, USA, AUS'
    $ws.Range('P28').Value = 'CAN, NLD, RUS, POL, NOR, DEU, USA, FRA, ITA'; $ws.Range('R28').Value = 'CAN, NLD, RUS, POL, NOR, DEU, USA, FRA, ITA'
    $ws.Range('P29').Value = 'SWE'; $ws.Range('R29').Value = 'SWE'
    $ws.Range('P31').Value = 'ARG, AUS, BEL, BRA, CAN, CHL, HRV, EGY, EST, DEU, GRC, HUN'; $ws.Range('R31').Value = 'ARG, AUS, BEL, BRA, CAN, CHL, HRV, EGY, EST, DEU, GRC, HUN'
    $ws.Range('P32').Value = 'ZAF'; $ws.Range('R32').Value = 'ZAF'
    $ws.Range('P33').Value = 'DNK'; $ws.Range('R33').Value = 'DNK'
    $ws.Range('P34').Value = 'KOR, VNM, KHM, JPN, CHN'; $ws.Range('R34').Value = 'KOR, VNM, KHM, JPN, CHN'
    $ws.Range('P35').Value = 'USA'; $ws.Range('R35').Value = 'USA'
    $ws.Range('P36').Value = 'GBR'; $ws.Range('R36').Value = 'GBR'
    $ws.Range('P38').Value = 'CAN, CHN'; $ws.Range('R38').Value = 'CAN, CHN'
    $ws.Range('P39').Value = 'BGD, JPN, TWN, KOR, CHN, IND, SGP, IRN, MNG, SGP, MYS, USA'; $ws.Range('R39').Value = 'BGD, JPN, TWN, KOR, CHN, IND, SGP, IRN, MNG, SGP, MYS, USA'
    $ws.Range('P40').Value = 'GBR, NLD, DEU, FRA, DNK, GRC, FIN, NOR, SWE, USA'; $ws.Range('R40').Value = 'GBR, NLD, DEU, FRA, DNK, GRC, FIN, NOR, SWE, USA'
    $ws.Range('P41').Value = 'USA'; $ws.Range('R41').Value = 'USA'
    $ws.Range('P42').Value = 'FIN'; $ws.Range('R42').Value = 'FIN'
    $ws.Range('P45').Value = 'JPN'; $ws.Range('R45').Value = 'JPN'
    $ws.Range('P47').Value = 'CHN, GHA, IND, MEX, RUS, ZAF, GBR, USA, AUT, BEL, DNK, FRA, DEU, GRC, ITA, CHE, NLD, ESP, SWE'; $ws.Range('R47').Value = 'CHN, GHA, IND, MEX, RUS, ZAF, GBR, USA, AUT, BEL, DNK, FRA, DEU, GRC, ITA, CHE, NLD, ESP, SWE'
    $ws.Range('P50').Value = 'BEL, CZE, DNK, FRO, FIN, FRA, DEU, GRC, IRL, ITA, LTU, NLD, NOR, POL, PRT, SVK, ESP, SWE, CHE, UKR, GBR'; $ws.Range('R50').Value = 'BEL, CZE, DNK, FRO, FIN, FRA, DEU, GRC, IRL, ITA, LTU, NLD, NOR, POL, PRT, SVK, ESP, SWE, CHE, UKR, GBR'
    $ws.Range('P51').Value = 'NOR, FIN, SWE, GBR, NLD, FRA, ESP, ITA'; $ws.Range('R51').Value = 'NOR, FIN, SWE, GBR, NLD, FRA, ESP, ITA'
    $ws.Range('P52').Value = 'ESP, SWE, GBR, DNK, NOR, FRA, DEU, NLD, GRC, ITA'; $ws.Range('R52').Value = 'ESP, SWE, GBR, DNK, NOR, FRA, DEU, NLD, GRC, ITA'
    $ws.Range('P53').Value = 'AUS'; $ws.Range('R53').Value = 'AUS'
    $ws.Range('P54').Value = 'USA, SWE, FIN'; $ws.Range('R54').Value = 'USA, SWE, FIN'
    $ws.Range('P56').Value = 'ITA, NLD, GBR, SWE'; $ws.Range('R56').Value = 'ITA, NLD, GBR, SWE'
    $ws.Range('P57').Value = 'CAN, FRA'; $ws.Range('R57').Value = 'CAN, FRA'
    $ws.Range('P59').Value = 'FIN, DNK, NLD, DEU, ESP, BEL, ITA, POL'; $ws.Range('R59').Value = 'FIN, DNK, NLD, DEU, ESP, BEL, ITA, POL'
}

# ---------------------------------------------------------------------------
# Step 2: The original commit also bumped the internal sheetId of the two
# renamed-content sheets from 6 -> 12 and 7 -> 13 (e.g. because they were
# recreated/copied a number of times upstream). Replicate that by copying
# each sheet forward (which mints a fresh sheetId = current max + 1) enough
# times to land on the same ids, then drop the original + helper copies and
# restore the names/position.
# ---------------------------------------------------------------------------

# Burn sheetIds 8,9,10,11 with throwaway sheets so the *next* two sheets we
# mint (via Copy) land on 12 and 13.
for ($i = 0; $i -lt 4; $i++) {
    $tmp = $wb.Worksheets.Add()
}

# Copy "Table 1 updated" to the end of the tab strip -> new sheetId 12.
$wb.Worksheets.Item("Table 1 updated").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Copy "Table 1 complete" to the end of the tab strip -> new sheetId 13.
$wb.Worksheets.Item("Table 1 complete").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Remove the four helper sheets.
$wb.Worksheets.Item("Sheet1").Delete()
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$wb.Worksheets.Item("Sheet4").Delete()

# Remove the original (now superseded) sheets.
$wb.Worksheets.Item("Table 1 updated").Delete()
$wb.Worksheets.Item("Table 1 complete").Delete()

# Rename the copies back to the original names (they land right after
# "Table 1", preserving tab order).
$wb.Worksheets.Item("Table 1 updated (2)").Name = "Table 1 updated"
$wb.Worksheets.Item("Table 1 complete (2)").Name = "Table 1 complete"
